# Weekly refresh of the Fruta/Hortaliza data: the D, L, M, N, O, P, Q, R, S, T
# columns of each data row (2..26) get reassigned to the values another row
# previously held (a permutation of the existing weekly snapshots), while the
# descriptive columns A, B, C, E, F, G, H, I, J, K stay as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map: destination row -> source row (values to copy, read from the ORIGINAL
# workbook state before any writes happen).
$rowMap = @{
    2  = 26
    3  = 17
    4  = 4
    5  = 3
    6  = 8
    7  = 23
    8  = 24
    9  = 7
    10 = 18
    11 = 19
    12 = 25
    13 = 6
    14 = 5
    15 = 20
    16 = 13
    17 = 16
    18 = 10
    19 = 9
    20 = 14
    21 = 2
    22 = 12
    23 = 11
    24 = 21
    25 = 22
    26 = 15
}

# Columns (by index) that carry the data which rotates between rows.
# D=4, L=12, M=13, N=14, O=15, P=16, Q=17, R=18, S=19, T=20
$cols = @(4, 12, 13, 14, 15, 16, 17, 18, 19, 20)

# First, snapshot the current ("before") values for every row/col we need,
# so overwriting a row doesn't corrupt data still needed as a source later.
$snapshot = @{}
for ($r = 2; $r -le 26; $r++) {
    foreach ($c in $cols) {
        $snapshot["$r,$c"] = $ws.Cells.Item($r, $c).Value2
    }
}

# Now write each destination row using the snapshotted source row's values.
foreach ($destRow in $rowMap.Keys) {
    $srcRow = $rowMap[$destRow]
    foreach ($c in $cols) {
        $ws.Cells.Item($destRow, $c).Value = $snapshot["$srcRow,$c"]
    }
}
